$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "26.799.42", "  +2.16%  ")
    ,@("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.681.46", "  +1.52%  ")
    ,@("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "0.9956", "  -0.84%  ")
    ,@("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "227.52", "  +3.79%  ")
    ,@("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.5316", "  +1.58%  ")
    ,@("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "0.9963", "  -0.82%  ")
    ,@("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.2671", "  +0.96%  ")
    ,@("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.06396", "  +1.13%  ")
    ,@("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "20.97", "  +1.51%  ")
    ,@("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.07696", "  -1.00%  ")
    ,@("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "4.680", "  +2.65%  ")
    ,@("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.667.23", "  -2.06%  ")
    ,@("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "1.890.43", "  +0.34%  ")
    ,@("Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.5764", "  +2.20%  ")
    ,@("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0₅8266", "  +2.16%  ")
    ,@("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "66.43", "  +1.63%  ")
    ,@("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "26.634.11", "  +1.56%  ")
    ,@("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.000", "  -0.40%  ")
    ,@("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "4.704", "  -0.48%  ")
    ,@("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "199.49", "  +3.36%  ")
    ,@("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "10.71", "  +4.12%  ")
    ,@("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "6.063", "  +0.58%  ")
    ,@("BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "0.9972", "  -0.77%  ")
    ,@("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "145.64", "  +0.79%  ")
    ,@("Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.1221", "  +1.41%  ")
    ,@("Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "7.319", "  +0.67%  ")
    ,@("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.574", "  +4.73%  ")
    ,@("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "16.24", "  +1.55%  ")
    ,@("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.05620", "  +0.50%  ")
    ,@("PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.305", "  +2.12%  ")
    ,@("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "3.551", "  +1.75%  ")
    ,@("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "3.455", "  +2.19%  ")
    ,@("LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "1.591", "  -0.40%  ")
    ,@("ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "0.9602", "  +1.46%  ")
    ,@("MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.802", "  +0.23%  ")
    ,@("HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.412", "  +0.46%  ")
    ,@("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.5767", "  +0.37%  ")
    ,@("VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.01608", "  +0.46%  ")
    ,@("FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "5.838", "  -1.95%  ")
    ,@("Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "1.041.91", "  +0.49%  ")
    ,@("PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "0.9995", "  -0.47%  ")
    ,@("TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.8365", "  -1.47%  ")
    ,@("Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "102.24", "  -0.24%  ")
    ,@("RocketPoolETH", "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth", "1.810.75", "  +0.90%  ")
    ,@("Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "58.75", "  +0.72%  ")
    ,@("BabyDogeCoin", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge", "0.0₈106", "  +0.93%  ")
    ,@("Frax", "https://coinranking.com/coin/KfWtaeV1W+frax-frax", "1.014", "  +1.57%  ")
    ,@("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "8.085", "  +0.83%  ")
    ,@("Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.4321", "  -0.75%  ")
    ,@("Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.05192", "  -2.38%  ")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = "'" + $data[$i][2]
    $ws.Cells.Item($row, 4).Style = "Normal"
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Style = "Normal"
}
